$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Map of cell reference -> new text value (all target cells are text-typed)
$updates = @{
    'D2' = '303.55'
    'E2' = '5.82%'
    'G2' = '18'
    'D3' = '32.02'
    'E3' = '9.73%'
    'G3' = '18'
    'D4' = '5.253'
    'E4' = '1.09%'
    'G4' = '18'
    'D5' = '0.07463'
    'E5' = '6.80%'
    'G5' = '18'
    'D6' = '7.859'
    'G6' = '18'
    'D7' = '3.816'
    'E7' = '7.36%'
    'G7' = '18'
    'D8' = '1.516'
    'E8' = '7.00%'
    'G8' = '18'
    'D9' = '0.9200'
    'E9' = '1.93%'
    'G9' = '18'
    'D10' = '0.01763'
    'E10' = '2,617.60%'
    'G10' = '18'
    'D11' = '0.1697'
    'E11' = '5.43%'
    'G11' = '18'
    'D12' = '0.07881'
    'E12' = '4.20%'
    'G12' = '18'
    'D13' = '0.08040'
    'E13' = '3.55%'
    'G13' = '18'
    'D14' = '0.03038'
    'E14' = '3.84%'
    'G14' = '18'
    'D15' = '0.09909'
    'E15' = '9.94%'
    'G15' = '18'
    'D16' = '0.001492'
    'E16' = '-7.02%'
    'G16' = '18'
    'D17' = '0.04615'
    'E17' = '1.88%'
    'G17' = '18'
    'D18' = '0.006240'
    'E18' = '-1.46%'
    'G18' = '18'
    'E19' = '0.36%'
    'G19' = '18'
    'D20' = '2.230'
    'E20' = '0.01%'
    'G20' = '18'
    'D21' = '0.3325'
    'E21' = '2.95%'
    'G21' = '18'
    'D22' = '0.1328'
    'E22' = '-0.27%'
    'G22' = '18'
    'D23' = '4.501'
    'E23' = '12.10%'
    'G23' = '18'
    'D24' = '0.1624'
    'E24' = '1.70%'
    'G24' = '18'
    'D25' = '0.001222'
    'E25' = '1.14%'
    'G25' = '18'
    'D26' = '0.004444'
    'E26' = '4.58%'
    'G26' = '18'
    'D27' = '0.0001402'
    'E27' = '20.16%'
    'G27' = '18'
    'D28' = '0.0001751'
    'E28' = '5.54%'
    'G28' = '18'
    'G29' = '18'
    'G30' = '18'
    'G31' = '18'
    'G32' = '18'
    'G33' = '18'
    'G34' = '18'
    'G35' = '18'
    'G36' = '18'
    'G37' = '18'
    'G38' = '18'
    'G39' = '18'
    'D40' = '0.04513'
    'E40' = '4.13%'
    'G40' = '18'
    'D41' = '0.007184'
    'E41' = '3.46%'
    'G41' = '18'
    'E42' = '8.13%'
    'G42' = '18'
    'D43' = '0.002203'
    'E43' = '6.22%'
    'G43' = '18'
    'D44' = '0.01286'
    'E44' = '10.40%'
    'G44' = '18'
    'D45' = '0.00006164'
    'E45' = '5.43%'
    'G45' = '18'
    'D46' = '0.7097'
    'E46' = '-63.22%'
    'G46' = '18'
    'D47' = '0.01302'
    'E47' = '-0.21%'
    'G47' = '18'
    'G48' = '18'
    'G49' = '18'
    'G50' = '18'
    'G51' = '18'
}

foreach ($cellRef in $updates.Keys) {
    $cell = $ws.Range($cellRef)
    # Force text format so Excel stores these as strings, not numbers,
    # matching the original inline-string cell type in the sheet.
    $cell.NumberFormat = "@"
    $cell.Value = $updates[$cellRef]
}
